$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B20 formula: add " + 5" to the existing sum
$ws.Range("B20").Formula = "=2 + 2.5 + 5"

# Copy B20's number format (currency-style) onto the new B21/B22 cells
# so they match the rest of column B instead of picking up a plain style
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)

# Add value for B21 (row for "?" label)
$ws.Range("B21").Value = 0

# Add formula for B22
$ws.Range("B22").Formula = "=5 + 2"

# Add three new week rows below the existing data
$ws.Range("A23").Value = "week 04-10/04/2016"
$ws.Range("A24").Value = "week 11-17/04/2016"
$ws.Range("A25").Value = "week 18-24/04/2016"

# Update the selected cell to match the new active cell
$ws.Range("B23").Select()
